# Update Name of Algo
# Apply the updated KNN-imputed values to the corresponding cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C6").Value  = -12.502
$ws.Range("B7").Value  = 5.375
$ws.Range("A9").Value  = -21.743
$ws.Range("B12").Value = 5.662
$ws.Range("B14").Value = 6.101999999999999
$ws.Range("C15").Value = -13.178
$ws.Range("A18").Value = -21.991
$ws.Range("A20").Value = -20.658
$ws.Range("B26").Value = 5.962000000000001
$ws.Range("A27").Value = -21.272
$ws.Range("B27").Value = 6.058000000000001
$ws.Range("B29").Value = 5.934
$ws.Range("C33").Value = -11.229
$ws.Range("A35").Value = -20.359
$ws.Range("C35").Value = -12.82
$ws.Range("B37").Value = 8.270000000000001
$ws.Range("B38").Value = 5.547000000000001
$ws.Range("C38").Value = -12.423
$ws.Range("C43").Value = -12.528
$ws.Range("C44").Value = -11.72
$ws.Range("C47").Value = -11.637
$ws.Range("B51").Value = 5.816
$ws.Range("C51").Value = -11.76
$ws.Range("B52").Value = 5.938
$ws.Range("B55").Value = 6.393999999999999
$ws.Range("C57").Value = -13.667
$ws.Range("C63").Value = -12.05
$ws.Range("A69").Value = -21.701
$ws.Range("B69").Value = 5.934
$ws.Range("B70").Value = 6.08
$ws.Range("C70").Value = -11.114
$ws.Range("A76").Value = -20.738
$ws.Range("A78").Value = -20.609
$ws.Range("B81").Value = 6.154999999999999
$ws.Range("A82").Value = -21.909
$ws.Range("A83").Value = -20.762
$ws.Range("B83").Value = 7.492
$ws.Range("C88").Value = -12.713
$ws.Range("A93").Value = -21.909
$ws.Range("C99").Value = -12.526
$ws.Range("B102").Value = 7.189
